# 16102024 11:39 add new condition of rc_index
#
# Update the "REMOTE CONTROL LINK" table:
#   - S04 (row 35) programming details change from the old Curtain/Dimmer/Relay
#     text to a new Fan/Lamp/Mixed text.
#   - S05 (row 36) programming details change to add the Dimmer/Relay scene
#     lines that used to live on S04.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testing Programming Details")

$ws.Range("C35").Value = "1: DEVICE FAN1 - FAN`n2: DEVICE FAN2 - LAMP`n3: DEVICE FAN3 - WHOLE`n4: SCENE Mixed Type"

$ws.Range("C36").Value = "1: DEVICE CURTAIN_1 - CLOSE`n2: DEVICE CURTAIN_2 - OPEN`n3: DEVICE CURTAIN_3 - WHOLE`n4: SCENE Dimmer Type - TOGGLE`n5: SCENE Relay Type - MOMENTARY"

# Move the active selection from C35 to C36, matching the saved cursor position.
$ws.Range("C36").Select()
